$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric updates (refreshed daily COVID case counts) ---
$ws.Cells.Item(4,2).Value = 2697634
$ws.Cells.Item(4,3).Value = 15823
$ws.Cells.Item(4,4).Value = 1124054
$ws.Cells.Item(4,5).Value = 1444534
$ws.Cells.Item(4,7).Value = 263
$ws.Cells.Item(4,8).Value = 129046
$ws.Cells.Item(5,2).Value = 1383678
$ws.Cells.Item(5,3).Value = 13190
$ws.Cells.Item(5,5).Value = 567289
$ws.Cells.Item(5,7).Value = 542
$ws.Cells.Item(5,8).Value = 58927
$ws.Cells.Item(7,2).Value = 582482
$ws.Cells.Item(7,3).Value = 14946
$ws.Cells.Item(7,4).Value = 345070
$ws.Cells.Item(7,5).Value = 220090
$ws.Cells.Item(7,7).Value = 418
$ws.Cells.Item(7,8).Value = 17322
$ws.Cells.Item(8,2).Value = 312654
$ws.Cells.Item(8,3).Value = 689
$ws.Cells.Item(8,7).Value = 155
$ws.Cells.Item(8,8).Value = 43730
$ws.Cells.Item(11,2).Value = 279393
$ws.Cells.Item(11,3).Value = 3394
$ws.Cells.Item(11,4).Value = 241229
$ws.Cells.Item(11,5).Value = 32476
$ws.Cells.Item(11,7).Value = 113
$ws.Cells.Item(11,8).Value = 5688
$ws.Cells.Item(12,2).Value = 240578
$ws.Cells.Item(12,3).Value = 142
$ws.Cells.Item(12,4).Value = 190248
$ws.Cells.Item(12,5).Value = 15563
$ws.Cells.Item(12,7).Value = 23
$ws.Cells.Item(12,8).Value = 34767
$ws.Cells.Item(16,2).Value = 199906
$ws.Cells.Item(16,3).Value = 1293
$ws.Cells.Item(16,4).Value = 173111
$ws.Cells.Item(16,5).Value = 21664
$ws.Cells.Item(16,7).Value = 16
$ws.Cells.Item(16,8).Value = 5131
$ws.Cells.Item(17,2).Value = 195559
$ws.Cells.Item(17,3).Value = 167
$ws.Cells.Item(17,5).Value = 7418
$ws.Cells.Item(44,2).Value = 32568
$ws.Cells.Item(44,3).Value = 752
$ws.Cells.Item(44,4).Value = 17580
$ws.Cells.Item(44,5).Value = 14241
$ws.Cells.Item(44,7).Value = 14
$ws.Cells.Item(44,8).Value = 747
$ws.Cells.Item(45,2).Value = 32125
$ws.Cells.Item(45,3).Value = 601
$ws.Cells.Item(45,4).Value = 8928
$ws.Cells.Item(45,5).Value = 22126
$ws.Cells.Item(45,7).Value = 57
$ws.Cells.Item(45,8).Value = 1071
$ws.Cells.Item(63,2).Value = 13907
$ws.Cells.Item(63,3).Value = 336
$ws.Cells.Item(63,4).Value = 9897
$ws.Cells.Item(63,5).Value = 3098
$ws.Cells.Item(63,7).Value = 7
$ws.Cells.Item(63,8).Value = 912
$ws.Cells.Item(69,2).Value = 11895
$ws.Cells.Item(69,3).Value = 90
$ws.Cells.Item(69,4).Value = 7770
$ws.Cells.Item(69,5).Value = 3776
$ws.Cells.Item(72,2).Value = 8879
$ws.Cells.Item(72,3).Value = 17
$ws.Cells.Item(72,5).Value = 491
$ws.Cells.Item(72,7).Value = 1
$ws.Cells.Item(72,8).Value = 250
$ws.Cells.Item(80,4).Value = 2039
$ws.Cells.Item(80,5).Value = 4179
$ws.Cells.Item(80,7).Value = 4
$ws.Cells.Item(80,8).Value = 148
$ws.Cells.Item(92,2).Value = 4299
$ws.Cells.Item(92,3).Value = 43
$ws.Cells.Item(92,4).Value = 3998
$ws.Cells.Item(92,5).Value = 191
$ws.Cells.Item(123,2).Value = 1462
$ws.Cells.Item(123,3).Value = 12
$ws.Cells.Item(123,4).Value = 974
$ws.Cells.Item(123,5).Value = 428
$ws.Cells.Item(128,2).Value = 1174
$ws.Cells.Item(128,3).Value = 2
$ws.Cells.Item(128,4).Value = 1031
$ws.Cells.Item(129,2).Value = 1132
$ws.Cells.Item(129,3).Value = 4
$ws.Cells.Item(129,4).Value = 882
$ws.Cells.Item(129,5).Value = 241
$ws.Cells.Item(129,8).Value = 9
$ws.Cells.Item(130,4).Value = 432
$ws.Cells.Item(130,5).Value = 392
$ws.Cells.Item(130,8).Value = 304
$ws.Cells.Item(153,2).Value = 548
$ws.Cells.Item(153,3).Value = 47
$ws.Cells.Item(153,4).Value = 315
$ws.Cells.Item(153,5).Value = 221
$ws.Cells.Item(153,7).Value = 1
$ws.Cells.Item(153,8).Value = 12
$ws.Cells.Item(154,2).Value = 526
$ws.Cells.Item(154,3).Value = 4
$ws.Cells.Item(154,4).Value = 472
$ws.Cells.Item(154,5).Value = 52
$ws.Cells.Item(154,8).Value = 2
$ws.Cells.Item(155,2).Value = 509
$ws.Cells.Item(155,4).Value = 183
$ws.Cells.Item(155,5).Value = 305
$ws.Cells.Item(155,8).Value = 21
$ws.Cells.Item(156,4).Value = 217
$ws.Cells.Item(156,5).Value = 271
$ws.Cells.Item(156,8).Value = 13

# --- Country name re-ranking swaps (ties / overtakes after the data refresh) ---
$ws.Cells.Item(44,1).Value = "Republica Dominicana"
$ws.Cells.Item(45,1).Value = "Bolivia"

$ws.Cells.Item(129,1).Value = "Jordania"
$ws.Cells.Item(130,1).Value = "Yemen"

$ws.Cells.Item(153,1).Value = "Montenegro"
$ws.Cells.Item(154,1).Value = "Reunion"
$ws.Cells.Item(155,1).Value = "Tanzania"
$ws.Cells.Item(156,1).Value = "Surinam"

$ws.Cells.Item(203,1).Value = "Santa Lucia"
$ws.Cells.Item(204,1).Value = "Laos"

$ws.Cells.Item(209,1).Value = "Groenlandia"
$ws.Cells.Item(210,1).Value = "Islas Malvinas"

# --- Header timestamp update ---
$ws.Cells.Item(1,1).Value = "Datos actualizados a 30 de Junio de 2020 a las 18:40"
